# BP-359 Bank excel statements upload
# Reconciliation.xlsx: convert REF_NO column (B) values to text-formatted
# shared strings, apply text number format to columns A-C, apply a
# 2-decimal numeric format to column E, move the active selection, and
# flip the page orientation to portrait (adds pageSetup element).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-format columns A (TRN_REF_NO), B (BANK_CODES) and C (REF_NO)
$ws.Range("A2:C6").NumberFormat = "@"

# New bank-code reference values for column B (stored as text)
$ws.Range("B2").Value = "199999"
$ws.Range("B3").Value = "288888"
$ws.Range("B4").Value = "388888"
$ws.Range("B5").Value = "488888"
$ws.Range("B6").Value = "588888"

# Column E (TXN_AMT) gets an explicit 2-decimal numeric format
$ws.Range("E2:E6").NumberFormat = "0.00"

# Move the active selection
$ws.Range("C3").Select()

# Page setup (portrait orientation)
$ws.PageSetup.Orientation = 1
